# Apply crypto price/volume updates (GitHub Actions daily refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.510.26"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("E3").Value = "  +5.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5114"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4161"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08755"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.43%  "
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").Value = "2.015.12"
$ws.Range("E13").Value = "  +5.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.598"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.487"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001115"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06528"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.220"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.81%  "
$ws.Range("D23").Value = "30.571.35"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.231"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("D26").Value = "2.250.64"
$ws.Range("E26").Value = "  +5.29%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.434"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.141"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.099"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.835"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.368"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02524"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.493"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06658"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.120"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2198"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6679"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.234"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6192"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.200"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.666"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.268"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.44%  "
